# "Dokumentasi hasil uji coba.xlsx" - fill in the ARITHMETIC_N CODING results table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 4): "Size" -> "Size (byte)" for both the original-file size
# column (D) and the compressed-file size column (F). The other headers
# (No / Original file / Compressed File / Compression Ratio) keep their text.
$ws.Range("D4").Value = "Size (byte)"
$ws.Range("F4").Value = "Size (byte)"

# --- First data row (row 5): the actual berk.raw result.
$ws.Range("C5").Value = "berk.raw"
$ws.Range("D5").Value = 1151150
$ws.Range("E5").Formula = '=CONCATENATE(C5,".arith_n")'
$ws.Range("F5").Value = 2047514
$ws.Range("G5").Value = -0.77
$ws.Range("G5").NumberFormat = "0%"

# --- Remaining rows (6-24): only the "<name>.arith_n" helper formula is filled
# in (as a shared formula across the block), C/D/F/G stay blank for now.
$ws.Range("E6:E24").Formula = '=CONCATENATE(C6,".arith_n")'

# --- Match the centered alignment used throughout the table for the newly
# populated data cells (C5:G24) -- mirrors the existing B5:B24/header styling.
$ws.Range("C5:G24").HorizontalAlignment = -4108
$ws.Range("G5").HorizontalAlignment = -4108

# --- Column widths: D and F widen to fit the new "Size (byte)" header / values.
$ws.Columns.Item(4).ColumnWidth = 10.6
$ws.Columns.Item(6).ColumnWidth = 11.6

# --- Selection moved from H14 to L14.
$ws.Range("L14").Select() | Out-Null
